# Add full keyboard and numeric keypad
# - Typography sheet: set G7, and add new typography rows 8-11
#   (Typography_Numeric_40px, Display, Keyboard, Mode)
# - Translation sheet: add new translation rows 28-38 for the
#   decoder-config screen and the on-screen keyboard / numeric keypad

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Typography"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Typography")

# Existing row 7 gains a Wildcard Characters value in G7
$ws1.Cells.Item(7, 7).Value = '!”#*"%&()''$+-@_, .:;?/~±×÷•º`´{}©£€^®¥_=[]¡¢|\¿><'

# New row 8: Typography_Numeric_40px
$ws1.Cells.Item(8, 2).Value = "Typography_Numeric_40px"
$ws1.Cells.Item(8, 3).Value = "consola.ttf"
$ws1.Cells.Item(8, 4).Value = 40
$ws1.Cells.Item(8, 5).Value = 4
$ws1.Cells.Item(8, 6).Value = "?"
$ws1.Cells.Item(8, 8).Value = "0-9,A-F"
$ws1.Range("B8:J8").Style = "Normal"

# New row 9: Display
$ws1.Cells.Item(9, 2).Value = "Display"
$ws1.Cells.Item(9, 3).Value = "Asap-Regular.ttf"
$ws1.Cells.Item(9, 4).Value = 28
$ws1.Cells.Item(9, 5).Value = 4
$ws1.Cells.Item(9, 6).Value = "?"
$ws1.Cells.Item(9, 7).Value = '!”"#*%&()''$+-@_, .:;?/~±×÷•º`´{}©£€^®¥_=[]¡¢|\¿><'
$ws1.Cells.Item(9, 8).Value = "a-z,A-Z,0-9"
$ws1.Range("B9:J9").Style = "Normal"

# New row 10: Keyboard
$ws1.Cells.Item(10, 2).Value = "Keyboard"
$ws1.Cells.Item(10, 3).Value = "Asap-Regular.ttf"
$ws1.Cells.Item(10, 4).Value = 20
$ws1.Cells.Item(10, 5).Value = 4
$ws1.Cells.Item(10, 6).Value = "?"
$ws1.Cells.Item(10, 7).Value = '!”#*"%&()''$+-@_, .:;?/~±×÷•º`´{}©£€^®¥_=[]¡¢|\¿><'
$ws1.Cells.Item(10, 8).Value = "a-z,A-Z,0-9"
$ws1.Range("B10:J10").Style = "Normal"

# New row 11: Mode
$ws1.Cells.Item(11, 2).Value = "Mode"
$ws1.Cells.Item(11, 3).Value = "Asap-Regular.ttf"
$ws1.Cells.Item(11, 4).Value = 20
$ws1.Cells.Item(11, 5).Value = 4
$ws1.Cells.Item(11, 6).Value = "?"
$ws1.Range("B11:J11").Style = "Normal"

# ---------------------------------------------------------------
# Sheet 2: "Translation"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Translation")

$ws2.Cells.Item(28, 2).Value = "DccConfigAddress"
$ws2.Cells.Item(28, 3).Value = "Typography_40px"
$ws2.Cells.Item(28, 4).Value = "Left"
$ws2.Cells.Item(28, 5).Value = "Address"
$ws2.Cells.Item(28, 6).Value = "LTR"

$ws2.Cells.Item(29, 2).Value = "DccConfigName"
$ws2.Cells.Item(29, 3).Value = "Typography_40px"
$ws2.Cells.Item(29, 4).Value = "Left"
$ws2.Cells.Item(29, 5).Value = "Name"
$ws2.Cells.Item(29, 6).Value = "LTR"

$ws2.Cells.Item(30, 2).Value = "DccConfigDescription"
$ws2.Cells.Item(30, 3).Value = "Typography_40px"
$ws2.Cells.Item(30, 4).Value = "Left"
$ws2.Cells.Item(30, 5).Value = "Description"
$ws2.Cells.Item(30, 6).Value = "LTR"

$ws2.Cells.Item(31, 2).Value = "wildcardTextIdMediumNumeric"
$ws2.Cells.Item(31, 3).Value = "Typography_40px"
$ws2.Cells.Item(31, 4).Value = "Left"
$ws2.Cells.Item(31, 5).Value = "<value>"
$ws2.Cells.Item(31, 6).Value = "LTR"

$ws2.Cells.Item(32, 2).Value = "DccConfigDecoder"
$ws2.Cells.Item(32, 3).Value = "Typography_40px"
$ws2.Cells.Item(32, 4).Value = "Left"
$ws2.Cells.Item(32, 5).Value = "Decoder"
$ws2.Cells.Item(32, 6).Value = "LTR"

$ws2.Cells.Item(33, 2).Value = "DccConfigConfig"
$ws2.Cells.Item(33, 3).Value = "Typography_40px"
$ws2.Cells.Item(33, 4).Value = "Left"
$ws2.Cells.Item(33, 5).Value = "Config"
$ws2.Cells.Item(33, 6).Value = "LTR"

$ws2.Cells.Item(34, 2).Value = "DccConfigAllCVs"
$ws2.Cells.Item(34, 3).Value = "Typography_40px"
$ws2.Cells.Item(34, 4).Value = "Left"
$ws2.Cells.Item(34, 5).Value = "All CVs"
$ws2.Cells.Item(34, 6).Value = "LTR"

$ws2.Cells.Item(35, 2).Value = "wildcardTextIdNumeric"
$ws2.Cells.Item(35, 3).Value = "Typography_Numeric_40px"
$ws2.Cells.Item(35, 4).Value = "Left"
$ws2.Cells.Item(35, 5).Value = "<id>"
$ws2.Cells.Item(35, 6).Value = "LTR"

$ws2.Cells.Item(36, 2).Value = "EnteredText"
$ws2.Cells.Item(36, 3).Value = "Display"
$ws2.Cells.Item(36, 4).Value = "Left"
$ws2.Cells.Item(36, 5).Value = "<placeHolder>"
$ws2.Cells.Item(36, 6).Value = "LTR"

$ws2.Cells.Item(37, 2).Value = "NumMode"
$ws2.Cells.Item(37, 3).Value = "Mode"
$ws2.Cells.Item(37, 4).Value = "Left"
$ws2.Cells.Item(37, 5).Value = "ABC"
$ws2.Cells.Item(37, 6).Value = "LTR"

$ws2.Cells.Item(38, 2).Value = "AlphaMode"
$ws2.Cells.Item(38, 3).Value = "Mode"
$ws2.Cells.Item(38, 4).Value = "Left"
$ws2.Cells.Item(38, 5).NumberFormat = "@"
$ws2.Cells.Item(38, 5).Value = "123"
$ws2.Cells.Item(38, 6).Value = "LTR"

# Strip the default-column style (s="1") picked up from column B's
# column-level style so the new cells match the un-styled data rows,
# and normalize row 38's forced-text number format back off.
$ws2.Range("B28:F38").Style = "Normal"

Write-Output "Edit complete"
